$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("city")

# Store price / sale price (column G) is 10 for every data row (rows 2-49)
$ws.Range("G2:G49").Value = 10

# Update the active selection to reflect where the user ended up (G2)
$ws.Range("G2").Select()
